$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to coerce numeric-looking price strings into text
# (matching the source data, which stores prices as plain text) without
# leaving a custom number-format style behind on the destination cells.
$helper = $ws.Range("Z100")

$ws.Range("D2").Value = "29.299.30"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "1.839.52"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$helper.NumberFormat = "@"
$helper.Value = "238.87"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E5").Value = "  -0.70%  "
$helper.NumberFormat = "@"
$helper.Value = "0.6249"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +0.05%  "
$helper.NumberFormat = "@"
$helper.Value = "0.07369"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E8").Value = "  -1.26%  "
$helper.NumberFormat = "@"
$helper.Value = "0.2887"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E9").Value = "  -0.81%  "
$helper.NumberFormat = "@"
$helper.Value = "24.75"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E10").Value = "  -1.03%  "
$helper.NumberFormat = "@"
$helper.Value = "0.07716"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.838.08"
$ws.Range("E12").Value = "  -0.76%  "
$ws.Range("E13").Value = "  -1.41%  "
$helper.NumberFormat = "@"
$helper.Value = "0.6627"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E14").Value = "  -2.99%  "
$helper.NumberFormat = "@"
$helper.Value = "0.00001044"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E15").Value = "  +2.08%  "
$helper.NumberFormat = "@"
$helper.Value = "81.33"
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E16").Value = "  -1.55%  "
$helper.NumberFormat = "@"
$helper.Value = "6.236"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "29.310.75"
$ws.Range("E18").Value = "  -0.85%  "
$helper.NumberFormat = "@"
$helper.Value = "235.78"
$helper.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E19").Value = "  +2.47%  "
$helper.NumberFormat = "@"
$helper.Value = "12.21"
$helper.Copy()
$ws.Range("D20").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E20").Value = "  -1.32%  "
$ws.Range("E21").Value = "  +0.14%  "
$helper.NumberFormat = "@"
$helper.Value = "7.244"
$helper.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E22").Value = "  -3.57%  "
$helper.NumberFormat = "@"
$helper.Value = "1.001"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E23").Value = "  +0.00%  "
$helper.NumberFormat = "@"
$helper.Value = "157.40"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("E25").Value = "  -1.19%  "
$helper.NumberFormat = "@"
$helper.Value = "0.1334"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E26").Value = "  -2.23%  "
$helper.NumberFormat = "@"
$helper.Value = "17.26"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E27").Value = "  -1.74%  "
$helper.NumberFormat = "@"
$helper.Value = "0.07123"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E28").Value = "  +7.79%  "
$helper.NumberFormat = "@"
$helper.Value = "1.477"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("E33").Value = "  +0.79%  "
$helper.NumberFormat = "@"
$helper.Value = "1.787"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E34").Value = "  -3.45%  "
$helper.NumberFormat = "@"
$helper.Value = "0.6880"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E35").Value = "  -1.48%  "
$helper.NumberFormat = "@"
$helper.Value = "2.581"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E36").Value = "  +0.70%  "
$helper.NumberFormat = "@"
$helper.Value = "0.01823"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E37").Value = "  -2.49%  "
$helper.NumberFormat = "@"
$helper.Value = "2.783"
$helper.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("D39").Value = "1.234.01"
$ws.Range("E39").Value = "  -1.60%  "
$helper.NumberFormat = "@"
$helper.Value = "6.727"
$helper.Copy()
$ws.Range("D40").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E40").Value = "  -0.90%  "
$helper.NumberFormat = "@"
$helper.Value = "0.9449"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E41").Value = "  +1.18%  "
$helper.NumberFormat = "@"
$helper.Value = "1.001"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "1.996.76"
$ws.Range("E43").Value = "  -0.67%  "
$helper.NumberFormat = "@"
$helper.Value = "101.17"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E44").Value = "  -0.21%  "
$helper.NumberFormat = "@"
$helper.Value = "65.03"
$helper.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$helper.NumberFormat = "@"
$helper.Value = "0.00000000117"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E46").Value = "  +1.20%  "
$helper.NumberFormat = "@"
$helper.Value = "6.918"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("E48").Value = "  -2.73%  "
$helper.NumberFormat = "@"
$helper.Value = "8.866"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$helper.Clear()
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("E51").Value = "  -1.60%  "
